# Updates cryptos list values (price/volume columns) per the Fri Jun 21 2024
# GitHub Actions data refresh. Every write is prefixed with a literal
# apostrophe (Excel's text quote-prefix) so numeric-looking strings such as
# "586.26" or "7.13" are stored as text, not auto-converted to numbers -
# matching the original inline-string cells. The immediate Style reset to
# "Normal" clears the transient quotePrefix formatting flag that the
# apostrophe trick leaves behind, so the cell's style index is left exactly
# as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "64.066.06"
Set-TextValue "E2" "  -1.48%  "
Set-TextValue "D3" "3.515.91"
Set-TextValue "E3" "  -0.01%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.06%  "
Set-TextValue "D5" "586.26"
Set-TextValue "E5" "  -1.12%  "
Set-TextValue "D6" "133.43"
Set-TextValue "E6" "  -0.34%  "
Set-TextValue "D7" "3.514.76"
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "E9" "  -0.72%  "
Set-TextValue "E10" "  -0.53%  "
Set-TextValue "D11" "7.13"
Set-TextValue "E11" "  -0.28%  "
Set-TextValue "D12" "0.377"
Set-TextValue "E12" "  -2.37%  "
Set-TextValue "D13" "4.112.32"
Set-TextValue "E13" "  -0.13%  "
Set-TextValue "D14" "27.49"
Set-TextValue "E14" "  -0.42%  "
Set-TextValue "E16" "  -1.65%  "
Set-TextValue "D17" "3.508.00"
Set-TextValue "E17" "  -0.31%  "
Set-TextValue "D18" "64.085.06"
Set-TextValue "E18" "  -1.46%  "
Set-TextValue "D19" "9.81"
Set-TextValue "E19" "  -2.79%  "
Set-TextValue "D20" "13.90"
Set-TextValue "E20" "  -2.91%  "
Set-TextValue "D21" "5.61"
Set-TextValue "E21" "  -1.05%  "
Set-TextValue "D22" "383.29"
Set-TextValue "E22" "  -2.27%  "
Set-TextValue "E23" "  -1.17%  "
Set-TextValue "D24" "3.657.13"
Set-TextValue "E24" "  -0.06%  "
Set-TextValue "D25" "73.87"
Set-TextValue "E25" "  -0.94%  "
Set-TextValue "E26" "  -0.16%  "
Set-TextValue "E27" "  -0.52%  "
Set-TextValue "E28" "  +2.71%  "
Set-TextValue "E29" "  -1.11%  "
Set-TextValue "D30" "7.48"
Set-TextValue "E30" "  -2.27%  "
Set-TextValue "E31" "  +0.03%  "
Set-TextValue "E32" "  +1.06%  "
Set-TextValue "E33" "  -1.71%  "
Set-TextValue "D34" "3.528.09"
Set-TextValue "E34" "  +0.12%  "
Set-TextValue "E35" "  +0.03%  "
Set-TextValue "D36" "23.55"
Set-TextValue "E36" "  -2.18%  "
Set-TextValue "E37" "  -0.13%  "
Set-TextValue "D38" "5.36"
Set-TextValue "E38" "  +1.81%  "
Set-TextValue "E39" "  -0.26%  "
Set-TextValue "D40" "6.93"
Set-TextValue "E40" "  -0.39%  "
Set-TextValue "D41" "160.57"
Set-TextValue "E41" "  -4.36%  "
Set-TextValue "D42" "0.0784"
Set-TextValue "E42" "  -2.59%  "
Set-TextValue "D43" "26.77"
Set-TextValue "E43" "  +4.48%  "
Set-TextValue "D44" "0.812"
Set-TextValue "E44" "  -1.03%  "
Set-TextValue "D45" "0.999"
Set-TextValue "E45" "  -0.13%  "
Set-TextValue "B46" "ONDO"
Set-TextValue "C46" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D46" "1.21"
Set-TextValue "E46" "  -3.45%  "
Set-TextValue "B47" "OKB"
Set-TextValue "C47" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D47" "41.66"
Set-TextValue "E47" "  -2.96%  "
Set-TextValue "E48" "  -0.77%  "
Set-TextValue "E49" "  -2.99%  "
Set-TextValue "D50" "2.474.76"
Set-TextValue "E50" "  +2.26%  "
Set-TextValue "D51" "6.80"
Set-TextValue "E51" "  -1.44%  "
